$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells remain text (avoid Excel auto-converting numeric-looking strings)
$cells = @('D2', 'E2', 'D3', 'E3', 'D4', 'E4', 'D5', 'E5', 'D6', 'E6', 'D7', 'E7', 'D8', 'D9', 'E9', 'D10', 'E10', 'E11', 'D12', 'E12', 'D13', 'E13', 'D14', 'E14', 'D15', 'E15', 'D16', 'E16', 'E17', 'D18', 'E18', 'D19', 'E19', 'D20', 'E20', 'D21', 'E21', 'D22', 'E22', 'E23', 'D24', 'E24', 'D25', 'E25', 'B26', 'C26', 'D26', 'E26', 'B27', 'C27', 'D27', 'E27', 'D28', 'E28', 'D29', 'E29', 'B30', 'C30', 'D30', 'E30', 'B31', 'C31', 'D31', 'E31', 'B32', 'C32', 'D32', 'E32', 'D33', 'E33', 'E34', 'D35', 'E35', 'D36', 'E36', 'D37', 'E37', 'D38', 'E38', 'D39', 'E39', 'D40', 'E40', 'D41', 'E41', 'D42', 'E42', 'D43', 'E43', 'D44', 'E44', 'D45', 'E45', 'B46', 'C46', 'D46', 'E46', 'B47', 'C47', 'D47', 'E47', 'D48', 'E48', 'D49', 'E49', 'D50', 'E50', 'D51', 'E51')
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '58.442.17'
$ws.Range('E2').Value = '  +0.99%  '
$ws.Range('D3').Value = '2.485.07'
$ws.Range('E3').Value = '  +1.27%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '521.68'
$ws.Range('E5').Value = '  +0.72%  '
$ws.Range('D6').Value = '132.75'
$ws.Range('E6').Value = '  +0.87%  '
$ws.Range('D7').Value = '0.995'
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('D8').Value = '0.557'
$ws.Range('D9').Value = '2.520.44'
$ws.Range('E9').Value = '  +2.57%  '
$ws.Range('D10').Value = '0.0977'
$ws.Range('E10').Value = '  -0.46%  '
$ws.Range('E11').Value = '  -0.19%  '
$ws.Range('D12').Value = '5.16'
$ws.Range('E12').Value = '  -1.84%  '
$ws.Range('D13').Value = '0.332'
$ws.Range('E13').Value = '  -2.03%  '
$ws.Range('D14').Value = '2.926.15'
$ws.Range('E14').Value = '  +1.45%  '
$ws.Range('D15').Value = '58.334.57'
$ws.Range('E15').Value = '  +0.94%  '
$ws.Range('D16').Value = '22.17'
$ws.Range('E16').Value = '  +0.02%  '
$ws.Range('E17').Value = '  +0.13%  '
$ws.Range('D18').Value = '2.498.47'
$ws.Range('E18').Value = '  +1.88%  '
$ws.Range('D19').Value = '10.75'
$ws.Range('E19').Value = '  +0.97%  '
$ws.Range('D20').Value = '321.66'
$ws.Range('E20').Value = '  +1.14%  '
$ws.Range('D21').Value = '4.17'
$ws.Range('E21').Value = '  +0.53%  '
$ws.Range('D22').Value = '6.04'
$ws.Range('E22').Value = '  +5.47%  '
$ws.Range('E23').Value = '  -0.44%  '
$ws.Range('D24').Value = '63.93'
$ws.Range('E24').Value = '  -0.46%  '
$ws.Range('D25').Value = '0.402'
$ws.Range('E25').Value = '  -1.18%  '
$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D26').Value = '0.161'
$ws.Range('E26').Value = '  +1.10%  '
$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').Value = '0.990'
$ws.Range('E27').Value = '  -0.69%  '
$ws.Range('D28').Value = '7.38'
$ws.Range('E28').Value = '  +0.78%  '
$ws.Range('D29').Value = '0.0₃0757'
$ws.Range('E29').Value = '  +2.95%  '
$ws.Range('B30').Value = 'Fetch.AI'
$ws.Range('C30').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D30').Value = '1.20'
$ws.Range('E30').Value = '  +3.03%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '1.71'
$ws.Range('E31').Value = '  +1.53%  '
$ws.Range('B32').Value = 'Monero'
$ws.Range('C32').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D32').Value = '167.11'
$ws.Range('E32').Value = '  +0.79%  '
$ws.Range('D33').Value = '6.26'
$ws.Range('E33').Value = '  +0.92%  '
$ws.Range('E34').Value = '  -0.14%  '
$ws.Range('D35').Value = '0.993'
$ws.Range('E35').Value = '  -0.43%  '
$ws.Range('D36').Value = '18.10'
$ws.Range('E36').Value = '  +0.39%  '
$ws.Range('D37').Value = '1.27'
$ws.Range('E37').Value = '  -2.06%  '
$ws.Range('D38').Value = '3.96'
$ws.Range('E38').Value = '  +0.36%  '
$ws.Range('D39').Value = '36.84'
$ws.Range('E39').Value = '  +2.30%  '
$ws.Range('D40').Value = '1.46'
$ws.Range('E40').Value = '  -0.34%  '
$ws.Range('D41').Value = '0.781'
$ws.Range('E41').Value = '  -0.45%  '
$ws.Range('D42').Value = '278.53'
$ws.Range('E42').Value = '  +2.71%  '
$ws.Range('D43').Value = '5.10'
$ws.Range('E43').Value = '  +1.81%  '
$ws.Range('D44').Value = '3.44'
$ws.Range('E44').Value = '  +0.42%  '
$ws.Range('D45').Value = '0.599'
$ws.Range('E45').Value = '  +1.84%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '122.55'
$ws.Range('E46').Value = '  -1.45%  '
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D47').Value = '0.0919'
$ws.Range('E47').Value = '  +1.50%  '
$ws.Range('D48').Value = '0.0502'
$ws.Range('E48').Value = '  +3.41%  '
$ws.Range('D49').Value = '17.86'
$ws.Range('E49').Value = '  +1.96%  '
$ws.Range('D50').Value = '0.0213'
$ws.Range('E50').Value = '  +1.85%  '
$ws.Range('D51').Value = '16.96'
$ws.Range('E51').Value = '  +1.93%  '